$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.487.47"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "3.712.73"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.28"
$ws.Range("E5").Value = "  -3.46%  "

$ws.Range("E6").Value = "  -2.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "651.19"
$ws.Range("E7").Value = "  -3.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.428"
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  -6.47%  "

$ws.Range("D11").Value = "3.709.30"
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.35"
$ws.Range("E12").Value = "  -2.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000308"
$ws.Range("E13").Value = "  +14.91%  "

$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("E15").Value = "  +1.60%  "

$ws.Range("D16").Value = "4.401.18"
$ws.Range("E16").Value = "  +0.23%  "

$ws.Range("D17").Value = "96.377.24"
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.84"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").Value = "3.696.56"
$ws.Range("E19").Value = "  -1.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.12"
$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.64"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.505"
$ws.Range("E22").Value = "  -8.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "520.58"
$ws.Range("E23").Value = "  +0.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.40"
$ws.Range("E24").Value = "  -1.53%  "

$ws.Range("E25").Value = "  -0.95%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.90"
$ws.Range("E26").Value = "  -0.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.47"
$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.36"
$ws.Range("E28").Value = "  +2.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.180"
$ws.Range("E29").Value = "  +6.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.01"
$ws.Range("E30").Value = "  -3.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.13"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.86"
$ws.Range("E33").Value = "  +6.49%  "

$ws.Range("E34").Value = "  -0.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.29"
$ws.Range("E36").Value = "  -3.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "645.48"
$ws.Range("E37").Value = "  +4.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.589"
$ws.Range("E38").Value = "  -1.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.79"
$ws.Range("E39").Value = "  -0.40%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.04"
$ws.Range("E43").Value = "  +3.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.160"
$ws.Range("E44").Value = "  -1.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.967"
$ws.Range("E45").Value = "  -0.20%  "

$ws.Range("E46").Value = "  +0.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.433"
$ws.Range("E47").Value = "  +1.96%  "

$ws.Range("E48").Value = "  -1.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.57"
$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.49"
$ws.Range("E50").Value = "  -1.61%  "

$ws.Range("E51").Value = "  +1.49%  "

# Row 41 becomes EnergySwap (was Filecoin)
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "41.06"
$ws.Range("E41").Value = "  -3.20%  "

# Row 42 becomes Filecoin (was EnergySwap)
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.82"
$ws.Range("E42").Value = "  +11.04%  "
